$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "Grigoli et al 2020 WP"
$ws.Range("B21").Value = "Use dispersion of infl forecasts as proxy for extent of anchoring, and find that this goes up after a mon pol surprise. Rationalize using a model with RE and sticky info."

$ws.Range("A22").Value = "Slobodyan Wouters, 2012, 2017"
$ws.Range("B22").Value = "estimate medium-scale DSGE with adaptive learning, fits much better than RE"

$ws.Rows.Item(21).RowHeight = 30

$ws.Range("A23").Select()
